# Ticket 79 - Fix implicit sheet cloning case when the number of items in
# the collection is 1. This adds a second rendered copy of the
# "${dvs.name}$@l=0" implicit-clone sheet (named "...@l=1") right after
# "Static3", mirroring what JETT's implicit cloning now produces when the
# bean collection backing the sheet contains exactly one item.

$wb = $excel.ActiveWorkbook

# Source sheet to clone and the sheet after which the clone should land.
# Single-quote these literal names so PowerShell does not try to expand
# "${...}" as a variable reference.
$src   = $wb.Worksheets.Item('${dvs.name}$@l=0')
$after = $wb.Worksheets.Item('Static3')

# Copy() duplicates the sheet (data, styles, merged cells, page setup, ...)
# and places the new sheet immediately after $after.
$src.Copy($null, $after)

# The freshly copied sheet becomes the active sheet and is named
# "${dvs.name}$@l=0 (2)" by default - rename it to the expected tab name.
$clone = $wb.Worksheets.Item('${dvs.name}$@l=0 (2)')
$clone.Name = '${dvs.name}$@l=1'

# Restore the original active/selected tab (Static1) so the workbook's
# tab-selection state matches the pre-edit workbook aside from the new
# sheet being appended.
$wb.Worksheets.Item('Static1').Activate()
